# The "for_category" column (originally column D) is being removed.
# Columns E ("measurement") and F ("type_of") shift left by one, so the
# header row becomes: parent | title_ENG | title_AMH | measurement | type_of
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "measurement"
$ws.Range("E1").Value = "type_of"
$ws.Cells.Item(1, 6).Clear()

# View/selection tweaks that came along with the edit.
$excel.ActiveWindow.Zoom = 110
$ws.Range("D6").Select()
